$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 290.25
$ws.Range("I6").Value = 53.666668
$ws.Range("K6").Value = 161.000004
$ws.Range("M6").Value = -49.00000399999999
$ws.Range("H19").Value = 580.625
$ws.Range("I19").Value = 728.61536
$ws.Range("J19").Value = 479.3684
$ws.Range("K19").Value = 728.61536
$ws.Range("L19").Value = 479.3684
$ws.Range("M19").Value = -553.61536
$ws.Range("N19").Value = -829.3684000000001
$ws.Range("H43").Value = 1968.16
$ws.Range("I43").Value = 2075
$ws.Range("J43").Value = 1896.9333
$ws.Range("K43").Value = 2075
$ws.Range("L43").Value = 1896.9333
$ws.Range("M43").Value = -2006
$ws.Range("N43").Value = -2034.9333
$ws.Range("H95").Value = 34905.75
$ws.Range("J95").Value = 34905.75
$ws.Range("L95").Value = 34905.75
$ws.Range("N95").Value = -40397.75
$ws.Range("H109").Value = 37500
$ws.Range("J109").Value = 37500
$ws.Range("L109").Value = 37500
$ws.Range("N109").Value = -40274
$ws.Range("H110").Value = 15826.8
$ws.Range("J110").Value = 15826.8
$ws.Range("L110").Value = 15826.8
$ws.Range("N110").Value = -24006.8
$ws.Range("H116").Value = 3500
$ws.Range("I116").Value = 3700
$ws.Range("K116").Value = 3700
$ws.Range("M116").Value = -258
$ws.Range("H121").Value = 837.96295
$ws.Range("I121").Value = 650
$ws.Range("J121").Value = 880.6818
$ws.Range("K121").Value = 1950
$ws.Range("L121").Value = 2642.0454
$ws.Range("M121").Value = -203
$ws.Range("N121").Value = -6136.0454
$ws.Range("H127").Value = 1261.5
$ws.Range("I127").Value = 1229.4
$ws.Range("J127").Value = 1273.8462
$ws.Range("K127").Value = 3688.2
$ws.Range("L127").Value = 3821.5386
$ws.Range("M127").Value = 1271.8
$ws.Range("N127").Value = -13741.5386
$ws.Range("H129").Value = 856.7692
$ws.Range("I129").Value = 709.36365
$ws.Range("J129").Value = 914.6786
$ws.Range("K129").Value = 2128.09095
$ws.Range("L129").Value = 2744.0358
$ws.Range("M129").Value = 2871.90905
$ws.Range("N129").Value = -12744.0358
$ws.Range("H132").Value = 4198.579
$ws.Range("I132").Value = 4003.2
$ws.Range("K132").Value = 12009.6
$ws.Range("M132").Value = -9479.599999999999
$ws.Range("H135").Value = 690.7143
$ws.Range("I135").Value = 577.55554
$ws.Range("J135").Value = 894.4
$ws.Range("K135").Value = 5197.99986
$ws.Range("L135").Value = 8049.599999999999
$ws.Range("M135").Value = -2662.99986
$ws.Range("N135").Value = -13119.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 13560
$ws.Range("J37").Value = 13200
$ws.Range("L37").Value = 13200
$ws.Range("N37").Value = -13746
$ws.Range("H46").Value = 2380.3333
$ws.Range("J46").Value = 2380.3333
$ws.Range("L46").Value = 2380.3333
$ws.Range("N46").Value = -3018.3333
$ws.Range("H57").Value = 3600
$ws.Range("I57").Value = 3600
$ws.Range("K57").Value = 3600
$ws.Range("M57").Value = -3116
$ws.Range("H104").Value = 26500
$ws.Range("J104").Value = 26500
$ws.Range("L104").Value = 26500
$ws.Range("N104").Value = -33488

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 54353.316
$ws.Range("I20").Value = 745.9
$ws.Range("J20").Value = 113917.11
$ws.Range("K20").Value = 745.9
$ws.Range("L20").Value = 113917.11
$ws.Range("M20").Value = -498.9
$ws.Range("N20").Value = -114411.11
$ws.Range("H94").Value = 5315.773
$ws.Range("I94").Value = 610.5833
$ws.Range("J94").Value = 10962
$ws.Range("K94").Value = 610.5833
$ws.Range("L94").Value = 10962
$ws.Range("M94").Value = -159.5833
$ws.Range("N94").Value = -11864

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 1060
$ws.Range("I35").Value = 1060
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1060
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -766
$ws.Range("N35").ClearContents()
$ws.Range("H41").Value = 9362
$ws.Range("I41").Value = 2564.5
$ws.Range("J41").Value = 14800
$ws.Range("K41").Value = 2564.5
$ws.Range("L41").Value = 14800
$ws.Range("M41").Value = -2136.5
$ws.Range("N41").Value = -15656
$ws.Range("H132").Value = 10579.611
$ws.Range("I132").Value = 17225.715
$ws.Range("J132").Value = 6350.273
$ws.Range("K132").Value = 51677.145
$ws.Range("L132").Value = 19050.819
$ws.Range("M132").Value = -49147.145
$ws.Range("N132").Value = -24110.819
$ws.Range("H141").Value = 76963544
$ws.Range("I141").Value = 25000
$ws.Range("J141").Value = 83375090
$ws.Range("K141").Value = 25000
$ws.Range("L141").Value = 83375090
$ws.Range("M141").Value = -19820
$ws.Range("N141").Value = -83385450

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 489.38095
$ws.Range("I2").Value = 773.5
$ws.Range("J2").Value = 110.55556
$ws.Range("K2").Value = 4641
$ws.Range("L2").Value = 663.33336
$ws.Range("M2").Value = -4528
$ws.Range("N2").Value = -889.33336
$ws.Range("H33").Value = 166.75
$ws.Range("I33").Value = 33.333332
$ws.Range("J33").Value = 567
$ws.Range("K33").Value = 199.999992
$ws.Range("L33").Value = 3402
$ws.Range("M33").Value = 83.00000800000001
$ws.Range("N33").Value = -3968
$ws.Range("H39").Value = 3300
$ws.Range("I39").Value = 800
$ws.Range("K39").Value = 2400
$ws.Range("M39").Value = -2106
$ws.Range("H86").Value = 662.5714
$ws.Range("I86").Value = 399.33334
$ws.Range("J86").Value = 860
$ws.Range("K86").Value = 1198.00002
$ws.Range("L86").Value = 2580
$ws.Range("M86").Value = -12.00001999999995
$ws.Range("N86").Value = -4952
$ws.Range("H87").Value = 2188.0667
$ws.Range("I87").Value = 2188.0667
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 6564.2001
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -5316.2001
$ws.Range("N87").ClearContents()
$ws.Range("H89").Value = 662.5714
$ws.Range("I89").Value = 399.33334
$ws.Range("J89").Value = 860
$ws.Range("K89").Value = 3594.00006
$ws.Range("L89").Value = 7740
$ws.Range("M89").Value = 2333.99994
$ws.Range("N89").Value = -19596
$ws.Range("H90").Value = 2188.0667
$ws.Range("I90").Value = 2188.0667
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 19692.6003
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -13452.6003
$ws.Range("N90").ClearContents()
$ws.Range("H98").Value = 375.0909
$ws.Range("I98").Value = 288.6
$ws.Range("J98").Value = 560.4286
$ws.Range("K98").Value = 865.8000000000001
$ws.Range("L98").Value = 1681.2858
$ws.Range("M98").Value = 632.1999999999999
$ws.Range("N98").Value = -4677.2858
$ws.Range("H113").Value = 716516.7
$ws.Range("I113").Value = 467.33334
$ws.Range("J113").Value = 1313224.5
$ws.Range("K113").Value = 1402.00002
$ws.Range("L113").Value = 3939673.5
$ws.Range("M113").Value = 767.9999800000001
$ws.Range("N113").Value = -3944013.5
$ws.Range("H132").Value = 3889800.8
$ws.Range("J132").Value = 66247.35000000001
$ws.Range("L132").Value = 596226.15
$ws.Range("N132").Value = -601286.15
$ws.Range("H137").Value = 6007768
$ws.Range("I137").Value = 93590.836
$ws.Range("J137").Value = 20201792
$ws.Range("K137").Value = 280772.508
$ws.Range("L137").Value = 60605376
$ws.Range("M137").Value = -275672.508
$ws.Range("N137").Value = -60615576

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 28692.309
$ws.Range("I64").Value = 13000
$ws.Range("K64").Value = 13000
$ws.Range("M64").Value = -12752
$ws.Range("H67").Value = 28692.309
$ws.Range("I67").Value = 13000
$ws.Range("K67").Value = 13000
$ws.Range("M67").Value = -12142
$ws.Range("H70").Value = 3950
$ws.Range("I70").Value = 3933.3333
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 3933.3333
$ws.Range("L70").Value = 4000
$ws.Range("M70").Value = -3663.3333
$ws.Range("N70").Value = -4540
$ws.Range("H73").Value = 3950
$ws.Range("I73").Value = 3933.3333
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 3933.3333
$ws.Range("L73").Value = 4000
$ws.Range("M73").Value = -2997.3333
$ws.Range("N73").Value = -5872
$ws.Range("H109").Value = 30331.166
$ws.Range("J109").Value = 30331.166
$ws.Range("L109").Value = 30331.166
$ws.Range("N109").Value = -32411.166
$ws.Range("H141").Value = 44927.668
$ws.Range("J141").Value = 46168.625
$ws.Range("L141").Value = 46168.625
$ws.Range("N141").Value = -56528.625

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1750.2
$ws.Range("I7").Value = 1685.3334
$ws.Range("J7").Value = 1847.5
$ws.Range("K7").Value = 1685.3334
$ws.Range("L7").Value = 1847.5
$ws.Range("M7").Value = -1573.3334
$ws.Range("N7").Value = -2071.5
$ws.Range("H105").Value = 11376.923
$ws.Range("J105").Value = 11376.923
$ws.Range("L105").Value = 11376.923
$ws.Range("N105").Value = -18364.923
$ws.Range("H122").Value = 61221.65
$ws.Range("I122").Value = 85475.664
$ws.Range("J122").Value = 3012
$ws.Range("K122").Value = 256426.992
$ws.Range("L122").Value = 9036
$ws.Range("M122").Value = -253976.992
$ws.Range("N122").Value = -13936
$ws.Range("H126").Value = 1750.2
$ws.Range("I126").Value = 1685.3334
$ws.Range("J126").Value = 1847.5
$ws.Range("K126").Value = 5056.0002
$ws.Range("L126").Value = 5542.5
$ws.Range("M126").Value = -2586.0002
$ws.Range("N126").Value = -10482.5
$ws.Range("H136").Value = 5829.4287
$ws.Range("I136").Value = 1630.625
$ws.Range("J136").Value = 11427.833
$ws.Range("K136").Value = 4891.875
$ws.Range("L136").Value = 34283.499
$ws.Range("M136").Value = -2341.875
$ws.Range("N136").Value = -39383.499

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 23682
$ws.Range("J104").Value = 23682
$ws.Range("L104").Value = 23682
$ws.Range("N104").Value = -30670
$ws.Range("H132").Value = 4163.7446
$ws.Range("I132").Value = 4861.9707
$ws.Range("J132").Value = 2337.6155
$ws.Range("K132").Value = 14585.9121
$ws.Range("L132").Value = 7012.8465
$ws.Range("M132").Value = -12055.9121
$ws.Range("N132").Value = -12072.8465
$ws.Range("H136").Value = 5125.48
$ws.Range("I136").Value = 16172.429
$ws.Range("J136").Value = 829.44446
$ws.Range("K136").Value = 48517.287
$ws.Range("L136").Value = 2488.33338
$ws.Range("M136").Value = -45967.287
$ws.Range("N136").Value = -7588.33338
$ws.Range("H140").Value = 54021.145
$ws.Range("J140").Value = 54021.145
$ws.Range("L140").Value = 54021.145
$ws.Range("N140").Value = -64381.145
$ws.Range("H141").Value = 59905.418
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 59905.418
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 59905.418
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -70265.41800000001
